# Automatische test-sync: 2025-06-19 21:25:50
$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 16 with the incoming complaint e-mail ---
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A16").Value = "Klacht over levering"
$wsLogs.Range("B16").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C16").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$wsLogs.Range("D16").Value = "Klacht / Probleem"
$wsLogs.Range("F16").Value = "2025-06-19 21:25:27"
$wsLogs.Range("G16").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$fcCategorie = $wsLogs.Range("D2:D15").FormatConditions
for ($i = 1; $i -le $fcCategorie.Count; $i++) {
    $fcCategorie.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D16"))
}

$fcBeantwoord = $wsLogs.Range("G2:G15").FormatConditions
for ($i = 1; $i -le $fcBeantwoord.Count; $i++) {
    $fcBeantwoord.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G16"))
}

# --- Dashboard sheet: swap rows 5/6 categories and add new row 9 ---
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A5").Value = "Openingstijden / Locatie"
$wsDash.Range("A6").Value = "Factuur / Administratie"

$wsDash.Range("A9").Value = "Klacht / Probleem"
$wsDash.Range("B9").Value = 1

# Update the chart series ranges to include the new Dashboard row
$co = $wsDash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$9,Dashboard!`$B`$2:`$B`$9,1)"
